$wb = $excel.ActiveWorkbook

$wsGlobal = $wb.Worksheets.Item("GLOBAL RESULTS")
$wsGlobal.Range("C6").Value = 57358.00562883465
$wsGlobal.Range("C7").Value = 57922.16039073943
$wsGlobal.Range("C8").Value = 52129.94435166549
$wsGlobal.Range("C12").Value = 45963.09371145877
$wsGlobal.Range("C13").Value = 45963.09371145877
$wsGlobal.Range("C14").Value = 33093.09371145877
$wsGlobal.Range("C15").Value = 32364.006420458772
$wsGlobal.Range("C16").Value = 31513.796420458777
$wsGlobal.Range("C20").Value = 562489.8859000112
$wsGlobal.Range("C21").Value = 568022.3541958446
$wsGlobal.Range("C22").Value = 511220.11877626023
$wsGlobal.Range("C26").Value = 450743.97294547706
$wsGlobal.Range("C27").Value = 450743.97294547706
$wsGlobal.Range("C28").Value = 324532.3874454771
$wsGlobal.Range("C29").Value = 317382.48356319196
$wsGlobal.Range("C30").Value = 309044.771666692

$wsWing = $wb.Worksheets.Item("WING")
$wsWing.Range("C8").Value = 4286.0
$wsWing.Range("D8").Value = 33.55769530397933
$wsWing.Range("C12").Value = 5015.0
$wsWing.Range("D12").Value = 56.274344831884356
$wsWing.Range("C13").Value = 4085.7142857142853
$wsWing.Range("D13").Value = 27.31651508878769
